$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Player Info" sheet at the very front of the workbook
#    (Worksheets.Add() with no args inserts before the currently active
#    sheet, which defaults to the first one).
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ------------------------------------------------------------------
# 2. Populate "Player Info" with header + single data row
# ------------------------------------------------------------------
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Reuse the existing bold/bordered header style already present in the
# workbook (same look as every other sheet's header row) by copying the
# format from an existing header cell instead of rebuilding it from
# scratch - keeps the shared style table untouched.
$battingSheet.Range("A1").Copy() | Out-Null
$playerInfo.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data row - the ID is a numeric-looking value but must stay text, so a
# leading apostrophe forces Excel to store it as a string.
$playerInfo.Range("A2").Value = "'6454"
$playerInfo.Range("B2").Value = "Shahnawaz Dahani"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ------------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, values become bare
#    match codes (text) instead of full scorecard URLs.
# ------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "'4592"
$battingSheet.Range("D3").Value = "'4641"

# ------------------------------------------------------------------
# 4. "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE transform.
# ------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").Value = "'4592"
$bowlingSheet.Range("B3").Value = "'4641"
